# Adds a "goal" adjustment for b2 and a matching "goal-adjustment" row
# to the b5 section, per commit message: "I think I added a goal adjustment"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "goal" info near the top (rows 12-16) ---
$ws.Range("D12").Value = "goal"
$ws.Range("D13").Value = -3
$ws.Range("E13").Value = 2
$ws.Range("D15").Value = "goal adjustment for b2"
$ws.Range("D16").Formula = "=D13-B3"
$ws.Range("E16").Formula = "=E13-C3"

# --- Update the "For B2" table header/weights (rows 27-31) ---
$ws.Range("A27").Value = "For B2 and goal -3 5"
$ws.Range("E28").Value = 20
$ws.Range("E29").Value = 0.1
$ws.Range("E30").Value = 0.3

# --- Insert a new "goal" row before the old "total" row (old row 32) ---
$ws.Rows(32).Insert()

$ws.Range("A32").Value = "goal"
$ws.Range("C32").Formula = "=-3-B3"
$ws.Range("D32").Value = 5
$ws.Range("E32").Value = 0.2
$ws.Range("F32").Formula = "=C32*E32"
$ws.Range("G32").Formula = "=D32*E32"

# The "total" row is now row 33; widen its weighted sums to include the new goal row.
$ws.Range("F33").Formula = "=SUM(F28:F32)"
$ws.Range("G33").Formula = "=SUM(G28:G32)"

# --- Insert a new "goal-adjustment" row before the old "self-velocity" row
#     of the b5 section. After the previous insert, that row now sits at 42. ---
$ws.Rows(42).Insert()

$ws.Range("A42").Value = "goal-adjustment"
$ws.Range("C42").Formula = "=-3+6.5"
$ws.Range("D42").Formula = "=5-6"
$ws.Range("E42").Value = 1
$ws.Range("F42").Formula = "=C42*E42"
$ws.Range("G42").Formula = "=D42*E42"

# The b5 section's sum row is now row 44; widen it to include the new row.
$ws.Range("C44").Formula = "=SUM(C38:C43)"
$ws.Range("D44").Formula = "=SUM(D38:D43)"
$ws.Range("F44").Formula = "=SUM(F38:F43)"
$ws.Range("G44").Formula = "=SUM(G38:G43)"

# Match the selection cell recorded in the saved file.
$ws.Range("C43").Select()

# The saved workbook also carries a (portrait) page setup.
$ws.PageSetup.Orientation = 1
